# MedicationCatalog.pptx — R5 conformity corrections for ValueSet / ObservationDefinition
# examples: reposition several shapes/connectors around the
# ClinicalUseDefinition boxes and fix a mislabeled "ClinicalDefinition" ->
# "ClinicalUseDefinition" text run.
#
# NOTE on numeric literals below: PowerPoint's Shape.Left/Top/Width/Height
# are single-precision (float32) on the COM surface, and the host truncates
# (toward zero) the point value * 12700 when it re-serialises to EMU. The
# literals used here were chosen so that, after that float32 round-trip,
# they land exactly on the target EMU value from the target OOXML (rather
# than the "obvious" nearest decimal, which can be off by 1 EMU).

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    throw "Shape with id $id not found"
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 34 "ZoneTexte 33" -------------------------------------------------
$sh34 = Get-ShapeById $s.Shapes 34
$sh34.Left = 268.05853271484375   # 3404343 EMU (was 3529374)
$sh34.Top  = 411.825927734375     # 5230189 EMU (was 5242109)

# --- Shape 144 "Rectangle 143" ----------------------------------------------
$sh144 = Get-ShapeById $s.Shapes 144
$sh144.Left  = 218.55520629882812 # 2775651 EMU (was 2978150)
$sh144.Width = 101.22630310058594 # 1285574 EMU (was 1176593)
$para2 = $sh144.TextFrame.TextRange.Paragraphs(2)
$para2.Runs(1).Text = "profile of ClinicalUseDefinition"

# --- Shape 154 "Connecteur : en angle 153" ----------------------------------
$sh154 = Get-ShapeById $s.Shapes 154
$sh154.Left   = 219.59103393554688 # 2788806 EMU (was 2779682)
$sh154.Top    = 396.0234069824219  # 5029497 EMU (was 5186629)
$sh154.Height = 25.077402114868164 # 318483 EMU (was 4218)

# --- Shape 163 "ZoneTexte 162" -----------------------------------------------
$sh163 = Get-ShapeById $s.Shapes 163
$sh163.Left = 245.53465270996094  # 3118290 EMU (was 3274155)

# --- Shape 169 "ZoneTexte 168" -----------------------------------------------
$sh169 = Get-ShapeById $s.Shapes 169
$sh169.Left = 291.83026123046875  # 3706244 EMU (was 3550379)

# --- Shape 133 "Connecteur : en angle 132" ----------------------------------
$sh133 = Get-ShapeById $s.Shapes 133
$sh133.Left   = 228.2222137451172  # 2898422 EMU (was 2815294)
$sh133.Top    = 412.4694519042969  # 5238362 EMU (was 5155234)
$sh133.Height = 51.38433074951172  # 652581 EMU (was 818837)
$sh133.Adjustments.Item(1) = 0.33942  # was 0.33373

# --- Shape 146 "Rectangle 145" -----------------------------------------------
$sh146 = Get-ShapeById $s.Shapes 146
$sh146.Left = 290.4382019042969   # 3688565 EMU (was 3522309)

# --- Shape 153 "Connecteur : en angle 152" ----------------------------------
$sh153 = Get-ShapeById $s.Shapes 153
$sh153.Left   = 184.61151123046875 # 2344566 EMU (was 2261438)
$sh153.Top    = 456.0801696777344  # 5792218 EMU (was 5709090)
$sh153.Height = 116.44575500488281 # 1478861 EMU (was 1645117)
$sh153.Adjustments.Item(1) = 0.18561  # was 0.18277

# --- Shape 170 "Connecteur : en angle 169" ----------------------------------
$sh170 = Get-ShapeById $s.Shapes 170
$sh170.Left   = 182.9578094482422  # 2323564 EMU (was 2240436)
$sh170.Top    = 457.73388671875    # 5813220 EMU (was 5730092)
$sh170.Height = 179.5999298095703  # 2280919 EMU (was 2447175)
$sh170.Adjustments.Item(1) = 0.15415  # was 0.15189

# --- Shape 172 "ZoneTexte 171" -----------------------------------------------
$sh172 = Get-ShapeById $s.Shapes 172
$sh172.Left = 450.07781982421875  # 5715988 EMU (was 5830289)

# --- Shape 174 "ZoneTexte 173" -----------------------------------------------
$sh174 = Get-ShapeById $s.Shapes 174
$sh174.Left = 408.5465393066406   # 5188541 EMU (was 5115804)

# --- Shape 175 "ZoneTexte 174" -----------------------------------------------
$sh175 = Get-ShapeById $s.Shapes 175
$sh175.Left = 474.4989929199219   # 6026137 EMU (was 6150829)
